$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.908.30"
$ws.Range("E2").Value = "  +0.44%  "

$ws.Range("D3").Value = "1.642.59"
$ws.Range("E3").Value = "  +0.56%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.82%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5065"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.85%  "

$ws.Range("E7").Value = "  +0.49%  "

$ws.Range("E8").Value = "  +0.21%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06416"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07780"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.307"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.79%  "

$ws.Range("D13").Value = "1.626.92"
$ws.Range("E13").Value = "  -0.52%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5455"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.00%  "

$ws.Range("D15").Value = "0.0₅7903"
$ws.Range("E15").Value = "  -0.27%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.75%  "

$ws.Range("D17").Value = "25.993.31"

$ws.Range("E18").Value = "  +0.56%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "197.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.415"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.68%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.051"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.006"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.862"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.69%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.27%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1149"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.37%  "

$ws.Range("E27").Value = "  +3.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.243"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.29%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05056"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.86%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.271"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.200"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.540"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.63%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.370"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.89%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8949"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.601"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.79%  "

$ws.Range("D38").Value = "1.130.97"
$ws.Range("E38").Value = "  -3.58%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01565"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.54%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.008"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.685"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.75%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8166"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.48%  "

$ws.Range("D44").Value = "0.0₈123"
$ws.Range("E44").Value = "  +7.83%  "

$ws.Range("D45").Value = "1.780.31"
$ws.Range("E45").Value = "  +0.53%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4553"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.96%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.90%  "

$ws.Range("E48").Value = "  -0.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05097"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.008"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.64%  "

$ws.Range("E51").Value = "  +3.16%  "
